$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-15 from 45175 to 45183
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
